# "finestra incidenza 7gg centrata su ultimo g"
#
# Column C ("somma mobile 7gg.") and column D ("somma mobile 7gg. per
# 100mila abitanti") were computed using a 7-day window CENTERED on each
# row (3 days before .. 3 days after). This re-centers the window so the
# 7 days end ON the current row (6 days before .. current day), i.e. a
# trailing rolling window, matching the "incidenza 7gg" convention used
# elsewhere (rolling sum of the last 7 days, anchored on the last day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 184
$windowSize = 7

# Column D is simply column C scaled by a constant factor (100000 /
# population). Rather than recomputing that factor (and risking tiny
# floating point rounding drift vs. the values Excel already stored),
# recover it from the cells that are already populated by building a
# lookup table of "new positives count" -> "per 100k inhabitants" value
# from the existing (not-yet-touched) sheet contents.
$perHundredK = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cVal = $ws.Range("C$r").Value()
    if ($cVal -ne $null) {
        $dVal = $ws.Range("D$r").Value()
        if (-not $perHundredK.ContainsKey($cVal)) {
            $perHundredK[$cVal] = $dVal
        }
    }
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $windowStart = $r - $windowSize + 1

    if ($windowStart -lt $firstRow) {
        # not enough trailing history yet for a full 7-day window
        $ws.Range("C$r").ClearContents()
        $ws.Range("D$r").ClearContents()
    } else {
        $sum = 0
        for ($i = $windowStart; $i -le $r; $i++) {
            $sum = $sum + $ws.Cells.Item($i, 2).Value()
        }
        $ws.Range("C$r").Value = $sum
        $ws.Range("D$r").Value = $perHundredK[$sum]
    }
}
